# Update countries & provincias Spain
# - Reorder "Groenlandia" / "Islas Malvinas" shared strings (Groenlandia now
#   precedes Islas Malvinas in the list)
# - Refresh the "datos actualizados" timestamp
# - Refresh the COVID numbers for several countries (Estados Unidos,
#   Sudafrica, Egipto, Israel, Costa de Marfil, Yemen, Mozambique, Angola)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the order of Groenlandia / Islas Malvinas -----------------------
# Row 210 held "Islas Malvinas" and row 211 held "Groenlandia"; after the
# edit row 210 holds "Groenlandia" and row 211 holds "Islas Malvinas" (the
# countries' numeric stats are identical, so only the labels move).
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# --- Update the "last updated" timestamp -----------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Julio de 2020 a las 22:52"

# --- Estados Unidos (row 4) --------------------------------------------
$ws.Range("B4").Value = 4011242
$ws.Range("C4").Value = 49813
$ws.Range("D4").Value = 1867629
$ws.Range("E4").Value = 1998882
$ws.Range("G4").Value = 897
$ws.Range("H4").Value = 144731

# --- Sudafrica (row 8) ---------------------------------------------------
$ws.Range("B8").Value = 381798
$ws.Range("C8").Value = 8170
$ws.Range("D8").Value = 208144
$ws.Range("E8").Value = 168286
$ws.Range("G8").Value = 195
$ws.Range("H8").Value = 5368

# --- Egipto (row 28) -------------------------------------------------------
$ws.Range("B28").Value = 89078
$ws.Range("C28").Value = 676
$ws.Range("D28").Value = 29473
$ws.Range("E28").Value = 55206
$ws.Range("G28").Value = 47
$ws.Range("H28").Value = 4399

# --- Israel (row 43) ------------------------------------------------------
$ws.Range("B43").Value = 54042
$ws.Range("C43").Value = 2039
$ws.Range("D43").Value = 22743
$ws.Range("E43").Value = 30874
$ws.Range("G43").Value = 10
$ws.Range("H43").Value = 425

# --- Costa de Marfil (row 69) ----------------------------------------------
$ws.Range("B69").Value = 14531
$ws.Range("C69").Value = 219
$ws.Range("D69").Value = 8857
$ws.Range("E69").Value = 5581
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 93

# --- Yemen (row 132) ---------------------------------------------------
$ws.Range("B132").Value = 1629
$ws.Range("C132").Value = 10
$ws.Range("D132").Value = 741
$ws.Range("E132").Value = 432
$ws.Range("G132").Value = 9
$ws.Range("H132").Value = 456

# --- Mozambique (row 135) ---------------------------------------------------
$ws.Range("D135").Value = 506
$ws.Range("E135").Value = 1019

# --- Angola (row 152) --------------------------------------------------
$ws.Range("B152").Value = 779
$ws.Range("C152").Value = 30
$ws.Range("E152").Value = 528
$ws.Range("G152").Value = 1
$ws.Range("H152").Value = 30
